# Apply the "new Comp40 co-infection variable" edit to the derived-variables table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row right after the existing "Comp39" row (row 108) for the
#    new "Comp40 / coinfection_any" derived variable, then populate it.
[void]$ws.Rows.Item(108).Insert()

$ws.Range("A108").Value = "Comp40"
$ws.Range("B108").Value = "coinfection_any"
$ws.Range("C108").Value = "Complications"
$ws.Range("D108").Value = "Any co-infection within +/- 2 weeks of COVID-19 dx"
$ws.Range("E108").Value = "0 = No; 1 = Yes; 99 = Unknown"

# 2. Grow the table (Table1) so it covers the new row too.
$lo = $ws.ListObjects.Item(1)
[void]$lo.Resize($ws.Range("A1:E238"))

# 3. Re-sort the whole table by column A ascending (Variable #), which is how
#    the table was already configured (sortState/sortCondition on column A).
$sortRange = $ws.Range("A1:E238")
[void]$sortRange.Sort($ws.Range("A1"), 1, $null, $null, 1, $null, 1, 1)

# 4. Update the view / selection to match the saved state after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 95
$win.ScrollColumn = 1
[void]$ws.Range("A110").Select()
